$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.799179801305371
$ws.Range("C2").Value = 0.08293557258618023
$ws.Range("D2").Value = 0.02163023039262413
$ws.Range("E2").Value = 0.07654431705199727
$ws.Range("F2").Value = 5.863167701030648
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.2344677533487669
$ws.Range("K2").Value = 0.8031111687674866
$ws.Range("M2").Value = 0.2915320632618261

$ws.Range("B3").Value = 0.7828399897699114
$ws.Range("C3").Value = 0.08223725791580705
$ws.Range("D3").Value = 0.02002985914074173
$ws.Range("E3").Value = 0.0773216924432667
$ws.Range("F3").Value = 5.662669389073955
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.2304004207806258
$ws.Range("K3").Value = 0.7884657395181307
$ws.Range("M3").Value = 0.2902225149972395

$ws.Range("B4").Value = 0.7735445472802667
$ws.Range("C4").Value = 0.0819253811699383
$ws.Range("D4").Value = 0.01907601518173152
$ws.Range("E4").Value = 0.0778505462825354
$ws.Range("F4").Value = 5.540014274568676
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.2279588057310491
$ws.Range("K4").Value = 0.7802822408273755
$ws.Range("M4").Value = 0.2896888111423515

$ws.Range("B5").Value = 0.7699415970101313
$ws.Range("C5").Value = 0.08182756097004074
$ws.Range("D5").Value = 0.0186944011046819
$ws.Range("E5").Value = 0.07807903105058145
$ws.Range("F5").Value = 5.490141141982718
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.2269777854045287
$ws.Range("K5").Value = 0.7771503040802941
$ws.Range("M5").Value = 0.2895392314254472

$ws.Range("B6").Value = 0.7693544938100558
$ws.Range("C6").Value = 0.0818130824605845
$ws.Range("D6").Value = 0.01863145778633424
$ws.Range("E6").Value = 0.07811775463631321
$ws.Range("F6").Value = 5.481866270475791
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.2268157290217587
$ws.Range("K6").Value = 0.7766424893711843
$ws.Range("M6").Value = 0.2895184932018218

$ws.Range("B7").Value = 0.7734952079948414
$ws.Range("C7").Value = 0.08192394357089938
$ws.Range("D7").Value = 0.01907084011644855
$ws.Range("E7").Value = 0.07785357516713276
$ws.Range("F7").Value = 5.539341228362446
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.2279455189098272
$ws.Range("K7").Value = 0.7802391815194625
$ws.Range("M7").Value = 0.2896865189980851

$ws.Range("B8").Value = 0.7933925774970305
$ws.Range("C8").Value = 0.082670458681946
$ws.Range("D8").Value = 0.02107235336583813
$ws.Range("E8").Value = 0.07680166543637235
$ws.Range("F8").Value = 5.793939814886784
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.2330537298675708
$ws.Range("K8").Value = 0.7978932417926217
$ws.Range("M8").Value = 0.2910243682619509

$ws.Range("B9").Value = 0.8382822427379324
$ws.Range("C9").Value = 0.08506772311520194
$ws.Range("D9").Value = 0.02523291214817647
$ws.Range("E9").Value = 0.07514736684023759
$ws.Range("F9").Value = 6.296985270176094
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.2435165698313142
$ws.Range("K9").Value = 0.8389568695930905
$ws.Range("M9").Value = 0.2957975146957814

$ws.Range("B10").Value = 0.8748771173451075
$ws.Range("C10").Value = 0.08740670714408338
$ws.Range("D10").Value = 0.02844348372111227
$ws.Range("E10").Value = 0.07418042196780306
$ws.Range("F10").Value = 6.669182482728189
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.2514809082080092
$ws.Range("K10").Value = 0.8730971014588249
$ws.Range("M10").Value = 0.3006223303459592

$ws.Range("B11").Value = 0.8923181208078006
$ws.Range("C11").Value = 0.08859817342992926
$ws.Range("D11").Value = 0.02993967384770713
$ws.Range("E11").Value = 0.07379438130722704
$ws.Range("F11").Value = 6.839140080640902
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.2551655899447383
$ws.Range("K11").Value = 0.8895003311529592
$ws.Range("M11").Value = 0.3031052282367952

$ws.Range("B12").Value = 0.8990373010483665
$ws.Range("C12").Value = 0.08906783240850302
$ws.Range("D12").Value = 0.03051154927837274
$ws.Range("E12").Value = 0.07365592996841563
$ws.Range("F12").Value = 6.9035959455735
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.2565698372903853
$ws.Range("K12").Value = 0.8958380041117948
$ws.Range("M12").Value = 0.3040869891392006

$ws.Range("B13").Value = 0.8975850980903886
$ws.Range("C13").Value = 0.08896585878504482
$ws.Range("D13").Value = 0.03038814721517724
$ws.Range("E13").Value = 0.07368540407208002
$ws.Range("F13").Value = 6.889709858845833
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.2562670088159962
$ws.Range("K13").Value = 0.894467453028966
$ws.Range("M13").Value = 0.3038736998732858

$ws.Range("B14").Value = 0.8928686111120498
$ws.Range("C14").Value = 0.08863644139402993
$ws.Range("D14").Value = 0.02998661517648316
$ws.Range("E14").Value = 0.07378283586886347
$ws.Range("F14").Value = 6.844440943435984
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.2552809386681076
$ws.Range("K14").Value = 0.8900192038597368
$ws.Range("M14").Value = 0.30318516505492

$ws.Range("B15").Value = 0.8899945727056036
$ws.Range("C15").Value = 0.08843707440330206
$ws.Range("D15").Value = 0.02974136055867405
$ws.Range("E15").Value = 0.07384352265475869
$ws.Range("F15").Value = 6.816725130367558
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.2546781084031551
$ws.Range("K15").Value = 0.8873109673224917
$ws.Range("M15").Value = 0.3027688310856007

$ws.Range("B16").Value = 0.8737533228088239
$ws.Range("C16").Value = 0.08733141934855837
$ws.Range("D16").Value = 0.02834643669174852
$ws.Range("E16").Value = 0.07420673321645843
$ws.Range("F16").Value = 6.658088647926775
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.2512413524754749
$ws.Range("K16").Value = 0.8720427262133228
$ws.Range("M16").Value = 0.3004658737157868

$ws.Range("B17").Value = 0.8639935084186448
$ws.Range("C17").Value = 0.08668588425614132
$ws.Range("D17").Value = 0.02749995592417775
$ws.Range("E17").Value = 0.07444333288542992
$ws.Range("F17").Value = 6.560937960493135
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.2491488578227177
$ws.Range("K17").Value = 0.8629001139481147
$ws.Range("M17").Value = 0.2991269468429252

$ws.Range("B18").Value = 0.8584546077202901
$ws.Range("C18").Value = 0.0863265767043373
$ws.Range("D18").Value = 0.02701643123035069
$ws.Range("E18").Value = 0.07458448560309172
$ws.Range("F18").Value = 6.505119748026203
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.2479511100294047
$ws.Range("K18").Value = 0.8577236255983394
$ws.Range("M18").Value = 0.2983839382442142

$ws.Range("B19").Value = 0.8565920435279111
$ws.Range("C19").Value = 0.08620697564542468
$ws.Range("D19").Value = 0.0268532880741148
$ws.Range("E19").Value = 0.07463314792732056
$ws.Range("F19").Value = 6.486230891248709
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.2475465671853385
$ws.Range("K19").Value = 0.855985037258705
$ws.Range("M19").Value = 0.2981370203815317

$ws.Range("B20").Value = 0.8650247239918372
$ws.Range("C20").Value = 0.08675336087448215
$ws.Range("D20").Value = 0.02758971744988514
$ws.Range("E20").Value = 0.07441762209687219
$ws.Range("F20").Value = 6.571273547251621
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.2493710066812866
$ws.Range("K20").Value = 0.8638648590043942
$ws.Range("M20").Value = 0.2992666716171257

$ws.Range("B21").Value = 0.8942508420063007
$ws.Range("C21").Value = 0.08873269661368965
$ws.Range("D21").Value = 0.03010440965546479
$ws.Range("E21").Value = 0.07375400794184195
$ws.Range("F21").Value = 6.857734866146529
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.2555703282152422
$ws.Range("K21").Value = 0.891322334282961
$ws.Range("M21").Value = 0.3033862760384878

$ws.Range("B22").Value = 0.9140203154407516
$ws.Range("C22").Value = 0.09013407387348593
$ws.Range("D22").Value = 0.03177888230653991
$ws.Range("E22").Value = 0.0733653740793887
$ws.Range("F22").Value = 7.045518971013394
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.2596740768104695
$ws.Range("K22").Value = 0.91000284776851
$ws.Range("M22").Value = 0.3063208684494327

$ws.Range("B23").Value = 0.9034076365116732
$ws.Range("C23").Value = 0.08937622166061487
$ws.Range("D23").Value = 0.03088229421464916
$ws.Range("E23").Value = 0.07356867279166757
$ws.Range("F23").Value = 6.945241981487982
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.2574790348253231
$ws.Range("K23").Value = 0.8999651958800428
$ws.Range("M23").Value = 0.3047324204470314

$ws.Range("B24").Value = 0.8645582867411292
$ws.Range("C24").Value = 0.08672281789162639
$ws.Range("D24").Value = 0.02754912653343666
$ws.Range("E24").Value = 0.07442922996948198
$ws.Range("F24").Value = 6.566600719721407
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.2492705568204627
$ws.Range("K24").Value = 0.8634284494902147
$ws.Range("M24").Value = 0.2992034187064192

$ws.Range("B25").Value = 0.8255062694881303
$ws.Range("C25").Value = 0.08431838694237115
$ws.Range("D25").Value = 0.02408109951926463
$ws.Range("E25").Value = 0.07555122609507592
$ws.Range("F25").Value = 6.160462162629244
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.2406378098331317
$ws.Range("K25").Value = 0.827153873776723
$ws.Range("M25").Value = 0.2942753584471589
